# "Remove requirement for main contributor (creator)"
#
# The workbook is a data-entry template whose header row marks required
# columns in bold (Title, Creator, Date Issued). This change removes the
# "required" marking (bold) from the "Creator" header in B2, and clears
# the stray Title value in A4 that belonged to a row lacking the
# (previously required) Creator field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Un-bold the "Creator" header - it is no longer a required field.
$ws.Range("B2").Font.Bold = $false

# Clear the orphaned "Title" entry on the row that is missing the
# Creator field (the row this fixture is named after).
$ws.Range("A4").ClearContents()

# Leave the selection where the edit finished.
$ws.Range("B4").Select()
